# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New K (strikeouts) values replacing the old Strike# values in column G,
# for rows 2-30 (game log rows, most recent first).
$kValues = @{
    2  = 0
    3  = 5
    4  = 5
    5  = 4
    6  = 3
    7  = 5
    8  = 4
    9  = 2
    10 = 4
    11 = 2
    12 = 5
    13 = 5
    14 = 9
    15 = 9
    16 = 6
    17 = 6
    18 = 7
    19 = 8
    20 = 7
    21 = 9
    22 = 6
    23 = 9
    24 = 4
    25 = 5
    26 = 9
    27 = 7
    28 = 3
    29 = 4
    30 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
